$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C6").Value = 60000349
$ws.Range("C7").Value = 60000349
$ws.Range("D8").Value = 259
$ws.Range("D9").Value = 259
$ws.Range("C10").Value = 60000349
$ws.Range("D11").Value = 259
$ws.Range("D12").Value = 259
$ws.Range("D13").Value = 259
$ws.Range("D14").Value = 259
$ws.Range("C15").Value = 60000350
$ws.Range("D16").Value = 260
